$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country label order in the shared strings (Lituania before Eslovenia) ---
$ws.Range("A123").Value = "Lituania"
$ws.Range("A124").Value = "Eslovenia"

# --- Swap country label order in the shared strings (Groenlandia before Islas Malvinas) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Julio de 2020 a las 09:39"

# --- Row 6: India ---
$ws.Range("B6").Value2 = 1119307
$ws.Range("C6").Value2 = 1200
$ws.Range("D6").Value2 = 700646
$ws.Range("E6").Value2 = 391147
$ws.Range("G6").Value2 = 11
$ws.Range("H6").Value2 = 27514

# --- Row 46: Singapur ---
$ws.Range("B46").Value2 = 48035
$ws.Range("C46").Value2 = 123
$ws.Range("E46").Value2 = 3922

# --- Row 53: Armenia ---
$ws.Range("B53").Value2 = 34981
$ws.Range("C53").Value2 = 104
$ws.Range("D53").Value2 = 23502
$ws.Range("E53").Value2 = 10829
$ws.Range("G53").Value2 = 9
$ws.Range("H53").Value2 = 650

# --- Row 67: Uzbekistan ---
$ws.Range("B67").Value2 = 16966
$ws.Range("C67").Value2 = 359
$ws.Range("E67").Value2 = 7599
$ws.Range("G67").Value2 = 3
$ws.Range("H67").Value2 = 88

# --- Row 100: Hungria ---
$ws.Range("B100").Value2 = 4339
$ws.Range("C100").Value2 = 6
$ws.Range("D100").Value2 = 3232
$ws.Range("E100").Value2 = 511

# --- Row 123: now Lituania's data ---
$ws.Range("B123").Value2 = 1947
$ws.Range("C123").Value2 = 15
$ws.Range("D123").Value2 = 1601
$ws.Range("H123").Value2 = 80

# --- Row 124: now Eslovenia's data ---
$ws.Range("B124").Value2 = 1946
$ws.Range("D124").Value2 = 1568
$ws.Range("E124").Value2 = 266
$ws.Range("H124").Value2 = 112
